$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset styling across the affected block (rows 14-23) back to Normal so that
# pre-existing cell formatting (e.g. the old yellow separator row / left-aligned
# zero cell) doesn't leak into cells that no longer need it.
$ws.Range("A14:C23").Style = "Normal"

# Row 14: R2 + Square -> Shift + G -> Group
$ws.Range("A14").Value = "R2 + Square"
$ws.Range("B14").Value = "Shift + G"
$ws.Range("C14").Value = "Group"

# Row 15: R2 + Triangle -> Ctrl + Shift + G -> Ungroup
$ws.Range("A15").Value = "R2 + Triangle"
$ws.Range("B15").Value = "Ctrl + Shift + G"
$ws.Range("C15").Value = "Ungroup"

# Row 16: R2 + X -> Ctrl + E -> Chop
$ws.Range("A16").Value = "R2 + X"
$ws.Range("B16").Value = "Ctrl + E"
$ws.Range("C16").Value = "Chop"

# Row 17: R2 + Circle -> Ctrl + J -> Consolidate
$ws.Range("A17").Value = "R2 + Circle"
$ws.Range("B17").Value = "Ctrl + J"
$ws.Range("C17").Value = "Consolidate"

# Row 18: "The below are not coded yet" separator (moved here), with yellow highlight
$ws.Range("A18").Value = "The below are not coded yet"
$ws.Range("B18").Value = ""
$ws.Range("C18").Value = ""
$ws.Range("A18:C18").Interior.Color = 65535

# Row 19: R2 + Left -> Ctrl + Z  -> Undo
$ws.Range("A19").Value = "R2 + Left"
$ws.Range("B19").Value = "Ctrl + Z "
$ws.Range("C19").Value = "Undo"

# Row 20: R2 + Right -> Ctrl + Y -> Redo
$ws.Range("A20").Value = "R2 + Right"
$ws.Range("B20").Value = "Ctrl + Y"
$ws.Range("C20").Value = "Redo"

# Row 21: R2 + Up -> Ctrl + A -> Select All
$ws.Range("A21").Value = "R2 + Up"
$ws.Range("B21").Value = "Ctrl + A"
$ws.Range("C21").Value = "Select All"

# Row 22: R2 + Down -> 0 (left aligned) -> Deactivate
$ws.Range("A22").Value = "R2 + Down"
$ws.Range("B22").Value = 0
$ws.Range("B22").HorizontalAlignment = -4131
$ws.Range("C22").Value = "Deactivate"

# Row 23: R2 + L1 -> Shift  -> Shift (hold)
$ws.Range("A23").Value = "R2 + L1"
$ws.Range("B23").Value = "Shift "
$ws.Range("C23").Value = "Shift (hold)"

# Update selection to match the new active block
$ws.Range("A14:XFD17").Select()
